$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.668551087379456
$ws.Range("B1").Value = 2.058618307113647
$ws.Range("C1").Value = 2.883825063705444
$ws.Range("D1").Value = 4.512523174285889
$ws.Range("E1").Value = 0.7563785314559937
